# "Complete the leads assign" - add a new lead row (row 10) to the Permit
# Data sheet, duplicating the CMDA/PP/NHRB/S/0630/2024 lead (row 2) into
# columns A:M and R:T, and filling in the Architect details (N, O, P) that
# were previously missing for this lead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns A:M: same lead data as row 2 -------------------------------
$ws.Range("A2:M2").Copy()
$ws.Range("A10").PasteSpecial()

# Row 10 uses the same 15.75pt row height as the other data rows.
$ws.Rows("10:10").RowHeight = 15.75

# --- Architect details (previously blank for this lead) ----------------
$ws.Range("N10").Value = "9.P.Elanchezhiyan. B.Arch.,M.T.P"
$ws.Range("N10").WrapText = $true

$ws.Range("O10").Value = "CMDA Regn. No. RA/Gr.I/19/06/276, No. 14/S2, Thirumurthy Nagar, Madananda puram,, Chennai- 600 125."

$ws.Range("P10").Value = "elan@rspindia.net"
$ws.Range("P10").WrapText = $true

# --- Document links (View PDF / View Approved Plan / View Approval Letter)
# Add the hyperlinks first, then stamp the same "Link" look-and-feel already
# used by every other row's link cells (R2:T9) over them, so row 10 reuses
# the existing Link cell style instead of growing the style table.
$ws.Hyperlinks.Add($ws.Range("R10"), "https://cmdachennai.gov.in/pdfs/OnlinePPAApprovalDetails/PP-NHRB-S-0630-2024/PlanPermit.pdf", "", "View PDF", "View PDF")
$ws.Hyperlinks.Add($ws.Range("S10"), "https://cmdachennai.gov.in/pdfs/OnlinePPAApprovalDetails/PP-NHRB-S-0630-2024/ApprovedPlan.pdf", "", "View Approved Plan", "View Approved Plan")
$ws.Hyperlinks.Add($ws.Range("T10"), "https://cmdachennai.gov.in/pdfs/OnlinePPAApprovalDetails/PP-NHRB-S-0630-2024/ApprovalLetter.pdf", "", "View Approval Letter", "View Approval Letter")

$ws.Range("R10:T10").Style = "Link"

# --- Update the view so the newly entered columns are in focus ---------
$ws.Range("O15").Select()
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
